$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has columns:
#   A COUNTER, B STRATEGY, C BUKRS, D BUTXT, E XBLNR, F BLDAT, G BUDAT,
#   H BLART, I SGTXT, J LIFNR, K KUNNR, L HKONT, M WRBTR_S, N WRBTR_H,
#   O WAERS, P MWSKZ, Q CUSTOM1
# The target removes the (unused/example) STRATEGY column (B) and the
# CUSTOM1 column (Q), shifting everything else left.
#
# Delete the rightmost column first (Q / CUSTOM1) so the index of column B
# is unaffected by the first delete.
$ws.Columns(17).Delete()   # Q: CUSTOM1
$ws.Columns(2).Delete()    # B: STRATEGY

# Column deletion leaves the worksheet's AutoFilter pointed at its old
# (now wrong) range, so clear and re-apply it over the new, shifted extent.
$ws.AutoFilterMode = $false
$ws.Range("E1:O11").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name (which Excel derives
# from the AutoFilter range) in sync with the new range as well.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=test_youtuber!`$E`$1:`$O`$11"
    }
}

# Restore the active selection recorded in the saved workbook.
$ws.Range("B2:B6").Select()
